# Adds TAS-diagram classification rows for "TrachyBasalt" and "TrachyDacite"
# average compositions to the GEOROC average-composition worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("SIO2(WT%)","TIO2(WT%)","AL2O3(WT%)","CR2O3(WT%)","FEOT(WT%)","CAO(WT%)","MGO(WT%)","MNO(WT%)","NIO(WT%)","K2O(WT%)","NA2O(WT%)","P2O5(WT%)","H2O(WT%)")

# ---- TrachyBasalt block (rows 76-79) ----
$ws.Cells.Item(76, 1).Value = "TrachyBasalt"
$ws.Cells.Item(76, 1).Font.Bold = $true

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(77, $i + 1).Value = $headers[$i]
}
$ws.Cells.Item(77, 14).Value = "TOT"

$row78 = @(49.1599928240624, 2.00135221435978, 16.1793838780455, 0.0354620706087903, 9.85323331702572, 8.4665174809513, 5.48044428322085, 0.173956647485043, 0.0120906638116253, 2.45386227940377, 3.65144520286169, 0.654495412723449, 1.15671641791045)
for ($i = 0; $i -lt $row78.Length; $i++) {
    $ws.Cells.Item(78, $i + 1).Value = $row78[$i]
}
$ws.Cells.Item(78, 14).Formula = "=SUM(A78:M78)"

$row79 = @(50.1007668501417, 2.03965206088426, 16.4890084983115, 0.0361407077082532, 10.0417944913679, 8.62854109568068, 5.5853234611223, 0.177285653170692, 0.0123220426589075, 2.50082180407866, 3.72132285351144, 0.667020481363753, 0)
for ($i = 0; $i -lt $row79.Length; $i++) {
    $ws.Cells.Item(79, $i + 1).Value = $row79[$i]
}
$ws.Cells.Item(79, 14).Formula = "=SUM(A79:M79)"

# ---- TrachyDacite block (rows 81-84) ----
$ws.Cells.Item(81, 1).Value = "TrachyDacite"
$ws.Cells.Item(81, 1).Font.Bold = $true

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(82, $i + 1).Value = $headers[$i]
}
$ws.Cells.Item(82, 14).Value = "TOT"

$row83 = @(65.1370418540092, 0.736724639642887, 15.8511239820027, 0.00405695561269142, 4.09043837352238, 2.7907539381042, 1.15751397097257, 0.0952083098759671, 0.002, 4.27234897751844, 4.20337982066165, 0.22801308622974, 0.551666666666667)
for ($i = 0; $i -lt $row83.Length; $i++) {
    $ws.Cells.Item(83, $i + 1).Value = $row83[$i]
}
$ws.Cells.Item(83, 14).Formula = "=SUM(A83:M83)"

$row84 = @(66.0829506266567, 0.747423226496519, 16.0813112426478, 0.00411587001523502, 4.1498390068849, 2.83128078054617, 1.17432318717952, 0.0965909083633603, 0.00202904365152988, 4.33439128497698, 4.26442057004116, 0.231324252540095, 0)
for ($i = 0; $i -lt $row84.Length; $i++) {
    $ws.Cells.Item(84, $i + 1).Value = $row84[$i]
}
$ws.Cells.Item(84, 14).Formula = "=SUM(A84:M84)"

# ---- View bookkeeping (match Excel's post-edit scroll/selection state) ----
$ws.StandardWidth = 11.58984375
$ws.Cells.Item(77, 9).Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 49 | Out-Null
